$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data (date serial, B, C, D) to append starting at row 256
$data = @(
    @(256, 44330, 0, 1, 15.22997258604935),
    @(257, 44331, 0, 1, 15.22997258604935),
    @(258, 44332, 0, 1, 15.22997258604935),
    @(259, 44333, 0, 1, 15.22997258604935),
    @(260, 44334, 0, 1, 15.22997258604935),
    @(261, 44335, 0, 1, 15.22997258604935),
    @(262, 44336, 2, 2, 30.45994517209869),
    @(263, 44337, 0, 2, 30.45994517209869),
    @(264, 44338, 1, 3, 45.68991775814803),
    @(265, 44339, 0, 3, 45.68991775814803),
    @(266, 44340, 0, 3, 45.68991775814803),
    @(267, 44341, 0, 3, 45.68991775814803),
    @(268, 44342, 0, 3, 45.68991775814803),
    @(269, 44343, 0, 1, 15.22997258604935)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

# Copy the formatting (style index 2: bold/centered/bordered date format) from
# column A of the last pre-existing row down across all the newly added rows.
$ws.Cells.Item(255, 1).Copy()
$ws.Range($ws.Cells.Item(256, 1), $ws.Cells.Item(269, 1)).PasteSpecial(-4122)
